$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Approver")

# Update the Password column values (B2/B3): "Bingo@1234567" -> "Bingo@12345"
$ws.Range("B2").Value = "Bingo@12345"
$ws.Range("B3").Value = "Bingo@12345"

# Swap the hyperlink targets for A2 and B3 so that:
#   A2 -> mailto:Bingo@1234567 (was mailto:jmuller@hl.com.test)
#   B3 -> mailto:jmuller@hl.com.test (was mailto:Bingo@1234567)
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.Address = "mailto:Bingo@1234567"
    }
    elseif ($addr -eq '$B$3') {
        $hl.Address = "mailto:jmuller@hl.com.test"
    }
}

# Move the active sheet / selection from EventExp!N3 to Approver!A3
[void]$ws.Range("A3").Select()
